$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in the ENERGY_CARRIERS database
$ws.Range("A17").Value = "Sunlight"
$ws.Range("A18").Value = "Ultraviolet"

# Update the active selection (cosmetic - matches saved selection in file)
$ws.Range("A27").Select()
